$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45178 -> 45179, i.e. 2023-09-09 -> 2023-09-10) for every data row.
$range = $ws.Range("C2:C45")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
